$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing cell values in columns F and G (rows 306-508) ---
$ws.Cells.Item(306, 6).Value2 = 77293
$ws.Cells.Item(336, 6).Value2 = 82018
$ws.Cells.Item(352, 6).Value2 = 307956
$ws.Cells.Item(352, 7).Value2 = 3574
$ws.Cells.Item(359, 6).Value2 = 321005
$ws.Cells.Item(364, 6).Value2 = 168772
$ws.Cells.Item(366, 6).Value2 = 339614
$ws.Cells.Item(366, 7).Value2 = 2843
$ws.Cells.Item(367, 6).Value2 = 769647
$ws.Cells.Item(367, 7).Value2 = 3924
$ws.Cells.Item(369, 6).Value2 = 235658
$ws.Cells.Item(369, 7).Value2 = 2611
$ws.Cells.Item(370, 6).Value2 = 180824
$ws.Cells.Item(370, 7).Value2 = 2044
$ws.Cells.Item(371, 6).Value2 = 160463
$ws.Cells.Item(371, 7).Value2 = 1973
$ws.Cells.Item(372, 6).Value2 = 178898
$ws.Cells.Item(372, 7).Value2 = 1861
$ws.Cells.Item(373, 6).Value2 = 351093
$ws.Cells.Item(373, 7).Value2 = 2389
$ws.Cells.Item(374, 6).Value2 = 775230
$ws.Cells.Item(374, 7).Value2 = 3435
$ws.Cells.Item(376, 6).Value2 = 223088
$ws.Cells.Item(376, 7).Value2 = 2233
$ws.Cells.Item(377, 6).Value2 = 177043
$ws.Cells.Item(377, 7).Value2 = 1829
$ws.Cells.Item(378, 6).Value2 = 157668
$ws.Cells.Item(378, 7).Value2 = 1551
$ws.Cells.Item(379, 6).Value2 = 181073
$ws.Cells.Item(379, 7).Value2 = 1615
$ws.Cells.Item(380, 6).Value2 = 346018
$ws.Cells.Item(380, 7).Value2 = 2034
$ws.Cells.Item(381, 6).Value2 = 748433
$ws.Cells.Item(383, 6).Value2 = 222940
$ws.Cells.Item(383, 7).Value2 = 1771
$ws.Cells.Item(384, 6).Value2 = 172540
$ws.Cells.Item(385, 6).Value2 = 151268
$ws.Cells.Item(386, 6).Value2 = 183523
$ws.Cells.Item(386, 7).Value2 = 1364
$ws.Cells.Item(387, 6).Value2 = 351952
$ws.Cells.Item(388, 6).Value2 = 731489
$ws.Cells.Item(390, 6).Value2 = 220262
$ws.Cells.Item(391, 6).Value2 = 178264
$ws.Cells.Item(391, 7).Value2 = 1210
$ws.Cells.Item(392, 6).Value2 = 222306
$ws.Cells.Item(392, 7).Value2 = 1242
$ws.Cells.Item(393, 6).Value2 = 309596
$ws.Cells.Item(393, 7).Value2 = 1243
$ws.Cells.Item(395, 6).Value2 = 753791
$ws.Cells.Item(398, 6).Value2 = 300870
$ws.Cells.Item(398, 7).Value2 = 1479
$ws.Cells.Item(399, 6).Value2 = 200829
$ws.Cells.Item(400, 6).Value2 = 150048
$ws.Cells.Item(400, 7).Value2 = 806
$ws.Cells.Item(401, 6).Value2 = 273278
$ws.Cells.Item(402, 6).Value2 = 723339
$ws.Cells.Item(403, 6).Value2 = 353975
$ws.Cells.Item(403, 7).Value2 = 735
$ws.Cells.Item(404, 6).Value2 = 225128
$ws.Cells.Item(405, 6).Value2 = 175108
$ws.Cells.Item(406, 6).Value2 = 171781
$ws.Cells.Item(407, 6).Value2 = 158564
$ws.Cells.Item(408, 6).Value2 = 306113
$ws.Cells.Item(409, 6).Value2 = 709210
$ws.Cells.Item(410, 6).Value2 = 365105
$ws.Cells.Item(411, 6).Value2 = 225643
$ws.Cells.Item(412, 6).Value2 = 177047
$ws.Cells.Item(413, 6).Value2 = 149963
$ws.Cells.Item(414, 6).Value2 = 149296
$ws.Cells.Item(415, 6).Value2 = 309109
$ws.Cells.Item(416, 6).Value2 = 674342
$ws.Cells.Item(417, 6).Value2 = 344483
$ws.Cells.Item(418, 6).Value2 = 202734
$ws.Cells.Item(418, 7).Value2 = 702
$ws.Cells.Item(419, 6).Value2 = 149735
$ws.Cells.Item(420, 6).Value2 = 139369
$ws.Cells.Item(420, 7).Value2 = 502
$ws.Cells.Item(422, 6).Value2 = 299016
$ws.Cells.Item(424, 6).Value2 = 267211
$ws.Cells.Item(425, 6).Value2 = 138313
$ws.Cells.Item(427, 6).Value2 = 90694
$ws.Cells.Item(427, 7).Value2 = 372
$ws.Cells.Item(429, 6).Value2 = 178708
$ws.Cells.Item(429, 7).Value2 = 459
$ws.Cells.Item(430, 6).Value2 = 175758
$ws.Cells.Item(432, 6).Value2 = 122707
$ws.Cells.Item(433, 6).Value2 = 87279
$ws.Cells.Item(434, 6).Value2 = 79280
$ws.Cells.Item(435, 6).Value2 = 83441
$ws.Cells.Item(437, 6).Value2 = 167483
$ws.Cells.Item(439, 6).Value2 = 89311
$ws.Cells.Item(440, 6).Value2 = 73754
$ws.Cells.Item(441, 6).Value2 = 68402
$ws.Cells.Item(442, 6).Value2 = 70642
$ws.Cells.Item(443, 6).Value2 = 106937
$ws.Cells.Item(444, 6).Value2 = 104229
$ws.Cells.Item(445, 6).Value2 = 84527
$ws.Cells.Item(446, 6).Value2 = 86752
$ws.Cells.Item(447, 6).Value2 = 67052
$ws.Cells.Item(448, 6).Value2 = 61546
$ws.Cells.Item(451, 6).Value2 = 86652
$ws.Cells.Item(453, 6).Value2 = 70164
$ws.Cells.Item(456, 6).Value2 = 50475
$ws.Cells.Item(458, 6).Value2 = 70839
$ws.Cells.Item(459, 6).Value2 = 59830
$ws.Cells.Item(460, 6).Value2 = 58434
$ws.Cells.Item(462, 6).Value2 = 43613
$ws.Cells.Item(464, 6).Value2 = 73624
$ws.Cells.Item(465, 6).Value2 = 61864
$ws.Cells.Item(467, 6).Value2 = 52248
$ws.Cells.Item(468, 6).Value2 = 41864
$ws.Cells.Item(469, 6).Value2 = 41165
$ws.Cells.Item(470, 6).Value2 = 43542
$ws.Cells.Item(471, 6).Value2 = 66907
$ws.Cells.Item(472, 6).Value2 = 51874
$ws.Cells.Item(473, 6).Value2 = 39876
$ws.Cells.Item(474, 6).Value2 = 45300
$ws.Cells.Item(475, 6).Value2 = 36644
$ws.Cells.Item(479, 6).Value2 = 42620
$ws.Cells.Item(481, 6).Value2 = 41378
$ws.Cells.Item(482, 6).Value2 = 36209
$ws.Cells.Item(483, 6).Value2 = 65322
$ws.Cells.Item(485, 6).Value2 = 13930
$ws.Cells.Item(490, 6).Value2 = 10735
$ws.Cells.Item(492, 6).Value2 = 13918
$ws.Cells.Item(494, 6).Value2 = 6311
$ws.Cells.Item(495, 6).Value2 = 10240
$ws.Cells.Item(497, 6).Value2 = 7552
$ws.Cells.Item(499, 6).Value2 = 11115
$ws.Cells.Item(500, 6).Value2 = 7324
$ws.Cells.Item(501, 6).Value2 = 5632
$ws.Cells.Item(502, 6).Value2 = 9950
$ws.Cells.Item(503, 6).Value2 = 7304
$ws.Cells.Item(504, 6).Value2 = 7161
$ws.Cells.Item(505, 6).Value2 = 8164
$ws.Cells.Item(506, 6).Value2 = 10419
$ws.Cells.Item(506, 7).Value2 = 10
$ws.Cells.Item(507, 6).Value2 = 6808
$ws.Cells.Item(507, 7).Value2 = 12
$ws.Cells.Item(508, 6).Value2 = 5540

# --- Append new rows 509-511 ---
# row 509
$ws.Cells.Item(509, 1).Value2 = 44403
$ws.Cells.Item(509, 2).Value2 = 392406
$ws.Cells.Item(509, 3).Value2 = 8694
$ws.Cells.Item(509, 4).Value2 = 51
$ws.Cells.Item(509, 5).Value2 = 12534
$ws.Cells.Item(509, 6).Value2 = 8767
$ws.Cells.Item(509, 7).Value2 = 21
$ws.Cells.Item(509, 1).NumberFormat = "yyyy-mm-dd"
# row 510
$ws.Cells.Item(510, 1).Value2 = 44404
$ws.Cells.Item(510, 2).Value2 = 392488
$ws.Cells.Item(510, 3).Value2 = 7399
$ws.Cells.Item(510, 4).Value2 = 82
$ws.Cells.Item(510, 5).Value2 = 12534
$ws.Cells.Item(510, 6).Value2 = 7090
$ws.Cells.Item(510, 7).Value2 = 18
$ws.Cells.Item(510, 1).NumberFormat = "yyyy-mm-dd"
# row 511
$ws.Cells.Item(511, 1).Value2 = 44405
$ws.Cells.Item(511, 2).Value2 = 392537
$ws.Cells.Item(511, 3).Value2 = 6965
$ws.Cells.Item(511, 4).Value2 = 49
$ws.Cells.Item(511, 5).Value2 = 12536
$ws.Cells.Item(511, 6).Value2 = 4566
$ws.Cells.Item(511, 7).Value2 = 14
$ws.Cells.Item(511, 1).NumberFormat = "yyyy-mm-dd"
